# Update the "dSF" column (F) values for a handful of rows.
# These correspond to a repull/recalculation of the data as noted in the
# commit message ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -1
$ws.Range("F3").Value  = 2
$ws.Range("F4").Value  = -7
$ws.Range("F11").Value = -9
$ws.Range("F13").Value = -3
$ws.Range("F22").Value = -2
$ws.Range("F24").Value = -2
$ws.Range("F26").Value = 5
$ws.Range("F29").Value = -4
$ws.Range("F33").Value = -1
$ws.Range("F34").Value = -8
